$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 200
$ws.Range("I21").Value = 200
$ws.Range("K21").Value = 200
$ws.Range("M21").Value = 268

$ws.Range("H23").Value = 200
$ws.Range("I23").Value = 200
$ws.Range("K23").Value = 200
$ws.Range("M23").Value = 34

$ws.Range("H34").Value = 1470
$ws.Range("I34").Value = 1462.5
$ws.Range("K34").Value = 1462.5
$ws.Range("M34").Value = -1259.5

$ws.Range("H36").Value = 1470
$ws.Range("I36").Value = 1462.5
$ws.Range("K36").Value = 1462.5
$ws.Range("M36").Value = -747.5

$ws.Range("H58").Value = 137
$ws.Range("I58").Value = 137
$ws.Range("K58").Value = 411
$ws.Range("M58").Value = -261

$ws.Range("H62").Value = 2700
$ws.Range("I62").Value = 2700
$ws.Range("K62").Value = 2700
$ws.Range("M62").Value = -2076

$ws.Range("H65").Value = 2700
$ws.Range("I65").Value = 2700
$ws.Range("K65").Value = 13500
$ws.Range("M65").Value = -10380

$ws.Range("H74").Value = 4997.5
$ws.Range("I74").Value = 4997.5
$ws.Range("K74").Value = 4997.5
$ws.Range("M74").Value = -4061.5

$ws.Range("H77").Value = 4997.5
$ws.Range("I77").Value = 4997.5
$ws.Range("K77").Value = 24987.5
$ws.Range("M77").Value = -20307.5

$ws.Range("H135").Value = 1445.5
$ws.Range("I135").Value = 1669.5
$ws.Range("J135").Value = 997.5
$ws.Range("K135").Value = 15025.5
$ws.Range("L135").Value = 8977.5
$ws.Range("M135").Value = -12490.5
$ws.Range("N135").Value = -14047.5

$ws.Range("H137").Value = 4499.5713
$ws.Range("I137").Value = 3833
$ws.Range("K137").Value = 11499
$ws.Range("M137").Value = -8949

$ws.Range("H138").Value = 3171.7307
$ws.Range("I138").Value = 2412.7273
$ws.Range("J138").Value = 3728.3333
$ws.Range("K138").Value = 7238.1819
$ws.Range("L138").Value = 11184.9999
$ws.Range("M138").Value = -2098.1819
$ws.Range("N138").Value = -21464.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1690.6
$ws.Range("I45").Value = 1667.6666
$ws.Range("K45").Value = 1667.6666
$ws.Range("M45").Value = -1290.6666

$ws.Range("H63").Value = 1999
$ws.Range("I63").Value = 1999
$ws.Range("K63").Value = 1999
$ws.Range("M63").Value = -1313

$ws.Range("H66").Value = 1999
$ws.Range("I66").Value = 1999
$ws.Range("K66").Value = 9995
$ws.Range("M66").Value = -6563

$ws.Range("H74").Value = 2955.15
$ws.Range("I74").Value = 2620.7144
$ws.Range("K74").Value = 2620.7144
$ws.Range("M74").Value = -1746.7144

$ws.Range("H77").Value = 2955.15
$ws.Range("I77").Value = 2620.7144
$ws.Range("K77").Value = 13103.572
$ws.Range("M77").Value = -8735.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 4692.5
$ws.Range("I22").Value = 4692.5
$ws.Range("K22").Value = 4692.5
$ws.Range("M22").Value = -4519.5

$ws.Range("H139").Value = 99995
$ws.Range("J139").Value = 99995
$ws.Range("L139").Value = 99995
$ws.Range("N139").Value = -110275

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 1000
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3000
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -550
$ws.Range("N122").Value = -10900

$ws.Range("H141").Value = 200000
$ws.Range("J141").Value = 200000
$ws.Range("L141").Value = 200000
$ws.Range("N141").Value = -210360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2229.3635
$ws.Range("J12").Value = 3300.2856
$ws.Range("L12").Value = 9900.856800000001
$ws.Range("N12").Value = -10246.8568

$ws.Range("H80").Value = 5828.5713
$ws.Range("J80").Value = 6000
$ws.Range("L80").Value = 18000
$ws.Range("N80").Value = -19872

$ws.Range("H83").Value = 5828.5713
$ws.Range("J83").Value = 6000
$ws.Range("L83").Value = 54000
$ws.Range("N83").Value = -63360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5250
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730

$ws.Range("H73").Value = 5250
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064

$ws.Range("H80").Value = 5166.6665
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 5250
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 5250
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -7246

$ws.Range("H83").Value = 5166.6665
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 5250
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 26250
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -36234

$ws.Range("H130").Value = 99995
$ws.Range("J130").Value = 99995
$ws.Range("L130").Value = 99995
$ws.Range("N130").Value = -110035

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2106.879
$ws.Range("J22").Value = 2404.4546
$ws.Range("L22").Value = 2404.4546
$ws.Range("N22").Value = -2994.4546

$ws.Range("H27").Value = 2106.879
$ws.Range("J27").Value = 2404.4546
$ws.Range("L27").Value = 2404.4546
$ws.Range("N27").Value = -2618.4546

$ws.Range("H40").Value = 4111.3335
$ws.Range("I40").Value = 5667
$ws.Range("K40").Value = 5667
$ws.Range("M40").Value = -5531

$ws.Range("H122").Value = 4640.2
$ws.Range("I122").Value = 4640.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13920.6
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11470.6
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 2579.8
$ws.Range("I136").Value = 1525.75
$ws.Range("J136").Value = 3282.5
$ws.Range("K136").Value = 4577.25
$ws.Range("L136").Value = 9847.5
$ws.Range("M136").Value = -2027.25
$ws.Range("N136").Value = -14947.5

$ws.Range("H139").Value = 50000
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -44860
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 50001
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 50001
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 7521
$ws.Range("I132").Value = 6603.875
$ws.Range("J132").Value = 9966.666999999999
$ws.Range("K132").Value = 19811.625
$ws.Range("L132").Value = 29900.001
$ws.Range("M132").Value = -17281.625
$ws.Range("N132").Value = -34960.001

$ws.Range("H139").Value = 74997.5
$ws.Range("I139").Value = 50000
$ws.Range("J139").Value = 99995
$ws.Range("K139").Value = 50000
$ws.Range("L139").Value = 99995
$ws.Range("M139").Value = -44860
$ws.Range("N139").Value = -110275
